$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 898.2857
$ws.Range("I2").Value = 898.2857
$ws.Range("K2").Value = 898.2857
$ws.Range("M2").Value = -785.2857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3023.875
$ws.Range("J17").Value = 3023.875
$ws.Range("L17").Value = 9071.625
$ws.Range("N17").Value = -9407.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3240.1333
$ws.Range("I64").Value = 3200.5
$ws.Range("K64").Value = 3200.5
$ws.Range("M64").Value = -2952.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3240.1333
$ws.Range("I67").Value = 3200.5
$ws.Range("K67").Value = 3200.5
$ws.Range("M67").Value = -2342.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 646.3
$ws.Range("I80").Value = 680.75
$ws.Range("J80").Value = 637.6875
$ws.Range("K80").Value = 2042.25
$ws.Range("L80").Value = 1913.0625
$ws.Range("M80").Value = -1044.25
$ws.Range("N80").Value = -3909.0625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 646.3
$ws.Range("I83").Value = 680.75
$ws.Range("J83").Value = 637.6875
$ws.Range("K83").Value = 6126.75
$ws.Range("L83").Value = 5739.1875
$ws.Range("M83").Value = -1134.75
$ws.Range("N83").Value = -15723.1875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 759.7778
$ws.Range("I98").Value = 729.75
$ws.Range("K98").Value = 729.75
$ws.Range("M98").Value = 768.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 5416.6665
$ws.Range("I106").Value = 4125
$ws.Range("J106").Value = 8000
$ws.Range("K106").Value = 4125
$ws.Range("L106").Value = 8000
$ws.Range("M106").Value = -3494
$ws.Range("N106").Value = -9262

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 759.7778
$ws.Range("I122").Value = 729.75
$ws.Range("K122").Value = 2189.25
$ws.Range("M122").Value = 260.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 864.9474
$ws.Range("I132").Value = 864.9474
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2594.8422
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -64.84220000000005
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4232.1035
$ws.Range("J138").Value = 4363.9585
$ws.Range("L138").Value = 13091.8755
$ws.Range("N138").Value = -23371.8755

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 8499
$ws.Range("I141").Value = 7998.4
$ws.Range("K141").Value = 23995.2
$ws.Range("M141").Value = -18815.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2609.2222
$ws.Range("I45").Value = 2609.2222
$ws.Range("K45").Value = 2609.2222
$ws.Range("M45").Value = -2232.2222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2031.8
$ws.Range("I61").Value = 1586.6666
$ws.Range("K61").Value = 1586.6666
$ws.Range("M61").Value = -1374.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1200
$ws.Range("J110").Value = 1200
$ws.Range("L110").Value = 1200
$ws.Range("N110").Value = -5290

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5659.6
$ws.Range("I122").Value = 6237.25
$ws.Range("J122").Value = 3349
$ws.Range("K122").Value = 18711.75
$ws.Range("L122").Value = 10047
$ws.Range("M122").Value = -16261.75
$ws.Range("N122").Value = -14947

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2450.2222
$ws.Range("I132").Value = 1692.6842
$ws.Range("K132").Value = 5078.0526
$ws.Range("M132").Value = -2548.0526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2031.8
$ws.Range("I136").Value = 1586.6666
$ws.Range("K136").Value = 4759.9998
$ws.Range("M136").Value = -2209.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3301
$ws.Range("I86").Value = 3301
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3301
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2178
$ws.Range("N86").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3301
$ws.Range("I89").Value = 3301
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 16505
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -10889
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 523.3333
$ws.Range("I94").Value = 523.3333
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 523.3333
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -72.33330000000001
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3049.375
$ws.Range("I31").Value = 2158.8
$ws.Range("K31").Value = 2158.8
$ws.Range("M31").Value = -1863.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3049.375
$ws.Range("I34").Value = 2158.8
$ws.Range("K34").Value = 2158.8
$ws.Range("M34").Value = -1956.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 21618.316
$ws.Range("I86").Value = 11001.625
$ws.Range("J86").Value = 29339.545
$ws.Range("K86").Value = 11001.625
$ws.Range("L86").Value = 29339.545
$ws.Range("M86").Value = -9878.625
$ws.Range("N86").Value = -31585.545

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 21618.316
$ws.Range("I89").Value = 11001.625
$ws.Range("J89").Value = 29339.545
$ws.Range("K89").Value = 55008.125
$ws.Range("L89").Value = 146697.725
$ws.Range("M89").Value = -49392.125
$ws.Range("N89").Value = -157929.725

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 3178.5715
$ws.Range("I18").Value = 750
$ws.Range("K18").Value = 2250
$ws.Range("M18").Value = -2081

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3519.6
$ws.Range("J34").Value = 6750
$ws.Range("L34").Value = 20250
$ws.Range("N34").Value = -20418

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 6738.8
$ws.Range("J55").Value = 14750
$ws.Range("L55").Value = 44250
$ws.Range("N55").Value = -44604

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 999.7778
$ws.Range("I140").Value = 999.7778
$ws.Range("K140").Value = 2999.3334
$ws.Range("M140").Value = 2180.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3127.5557
$ws.Range("I132").Value = 2441.75
$ws.Range("K132").Value = 7325.25
$ws.Range("M132").Value = -4795.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 94500
$ws.Range("J128").Value = 94500
$ws.Range("L128").Value = 94500
$ws.Range("N128").Value = -104460

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5034.2383
$ws.Range("I132").Value = 4885.25
$ws.Range("J132").Value = 5511
$ws.Range("K132").Value = 14655.75
$ws.Range("L132").Value = 16533
$ws.Range("M132").Value = -12125.75
$ws.Range("N132").Value = -21593

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5486.769
$ws.Range("I81").Value = 3030.3635
$ws.Range("K81").Value = 6060.727
$ws.Range("M81").Value = -4999.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 5486.769
$ws.Range("I84").Value = 3030.3635
$ws.Range("K84").Value = 30303.635
$ws.Range("M84").Value = -24999.635

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1392.1428
$ws.Range("I100").Value = 1392.1428
$ws.Range("K100").Value = 2784.2856
$ws.Range("M100").Value = -2243.2856
